$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Fgf15-Klb, MuSCs -> MuSCs)
$ws.Range("M2").Value = 0.2921476666666666
$ws.Range("N2").Value = 0.876443
$ws.Range("O2").Value = 0.1567859736747012
$ws.Range("P2").Value = 0.1567859736747012
$ws.Range("Q2").Value = 0.03913317995
$ws.Range("R2").Value = 0.35219861955
$ws.Range("S2").Value = 0.1567859736747012
$ws.Range("T2").Value = 0.1567859736747012

# Row 3 (Fgf15-Klb, MuSCs -> ECs)
$ws.Range("O3").Value = 0.6778553360786825
$ws.Range("P3").Value = 0.6778553360786825
$ws.Range("S3").Value = 0.6778553360786825
$ws.Range("T3").Value = 0.6778553360786825

# Row 4 (Fgf15-Klb, MuSCs -> MuSCs)
$ws.Range("O4").Value = 0.1653586902466163
$ws.Range("P4").Value = 0.1653586902466163
$ws.Range("S4").Value = 0.1653586902466163
$ws.Range("T4").Value = 0.1653586902466163
